$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 48.4
$ws.Range("I11").Value = 48.4
$ws.Range("K11").Value = 48.4
$ws.Range("M11").Value = 91.59999999999999
$ws.Range("H12").Value = 364.57144
$ws.Range("I12").Value = 190
$ws.Range("K12").Value = 190
$ws.Range("M12").Value = -20
$ws.Range("H17").Value = 2439986.8
$ws.Range("J17").Value = 2439986.8
$ws.Range("L17").Value = 7319960.399999999
$ws.Range("N17").Value = -7320296.399999999
$ws.Range("H86").Value = 3235.25
$ws.Range("I86").Value = 2779.625
$ws.Range("K86").Value = 2779.625
$ws.Range("M86").Value = -1656.625
$ws.Range("H89").Value = 3235.25
$ws.Range("I89").Value = 2779.625
$ws.Range("K89").Value = 13898.125
$ws.Range("M89").Value = -8282.125
$ws.Range("H106").Value = 1846.5
$ws.Range("I106").Value = 1758.125
$ws.Range("K106").Value = 1758.125
$ws.Range("M106").Value = -1127.125
$ws.Range("H112").Value = 3968.25
$ws.Range("J112").Value = 3968.25
$ws.Range("L112").Value = 11904.75
$ws.Range("N112").Value = -14120.75
$ws.Range("H133").Value = 73055
$ws.Range("J133").Value = 73055
$ws.Range("L133").Value = 73055
$ws.Range("N133").Value = -83175
$ws.Range("H141").Value = 2423.3
$ws.Range("I141").Value = 2174.12
$ws.Range("K141").Value = 6522.36
$ws.Range("M141").Value = -1342.36
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 2436.2307
$ws.Range("I5").Value = 3806.625
$ws.Range("J5").Value = 243.6
$ws.Range("K5").Value = 3806.625
$ws.Range("L5").Value = 243.6
$ws.Range("M5").Value = -3694.625
$ws.Range("N5").Value = -467.6
$ws.Range("H45").Value = 7799.2
$ws.Range("I45").Value = 3666.6667
$ws.Range("J45").Value = 9570.286
$ws.Range("K45").Value = 3666.6667
$ws.Range("L45").Value = 9570.286
$ws.Range("M45").Value = -3289.6667
$ws.Range("N45").Value = -10324.286
$ws.Range("H97").Value = 1014.2353
$ws.Range("I97").Value = 759.5333000000001
$ws.Range("J97").Value = 2924.5
$ws.Range("K97").Value = 759.5333000000001
$ws.Range("L97").Value = 2924.5
$ws.Range("M97").Value = -263.5333000000001
$ws.Range("N97").Value = -3916.5
$ws.Range("H122").Value = 2241.7856
$ws.Range("I122").Value = 1656.8889
$ws.Range("J122").Value = 3294.6
$ws.Range("K122").Value = 4970.6667
$ws.Range("L122").Value = 9883.799999999999
$ws.Range("M122").Value = -2520.6667
$ws.Range("N122").Value = -14783.8
$ws.Range("H132").Value = 5573.6943
$ws.Range("I132").Value = 4707.9375
$ws.Range("K132").Value = 14123.8125
$ws.Range("M132").Value = -11593.8125
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 2436.2307
$ws.Range("I4").Value = 3806.625
$ws.Range("J4").Value = 243.6
$ws.Range("K4").Value = 3806.625
$ws.Range("L4").Value = 243.6
$ws.Range("M4").Value = -3691.625
$ws.Range("N4").Value = -473.6
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 261.91666
$ws.Range("I22").Value = 268.1
$ws.Range("J22").Value = 231
$ws.Range("K22").Value = 268.1
$ws.Range("L22").Value = 231
$ws.Range("M22").Value = 81.89999999999998
$ws.Range("N22").Value = -931
$ws.Range("H94").Value = 925.25
$ws.Range("I94").Value = 841.5
$ws.Range("K94").Value = 841.5
$ws.Range("M94").Value = -390.5
$ws.Range("H132").Value = 2671.25
$ws.Range("I132").Value = 2408.6
$ws.Range("J132").Value = 3984.5
$ws.Range("K132").Value = 7225.799999999999
$ws.Range("L132").Value = 11953.5
$ws.Range("M132").Value = -4695.799999999999
$ws.Range("N132").Value = -17013.5
$ws.Range("H134").Value = 6642.5264
$ws.Range("I134").Value = 5075.5625
$ws.Range("J134").Value = 14999.667
$ws.Range("K134").Value = 15226.6875
$ws.Range("L134").Value = 44999.001
$ws.Range("M134").Value = -12691.6875
$ws.Range("N134").Value = -50069.001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 994.4286
$ws.Range("I5").Value = 491.66666
$ws.Range("K5").Value = 1474.99998
$ws.Range("M5").Value = -1362.99998
$ws.Range("H23").Value = 7226.143
$ws.Range("I23").Value = 85.833336
$ws.Range("K23").Value = 257.500008
$ws.Range("M23").Value = -22.50000799999998
$ws.Range("H126").Value = 7663.5
$ws.Range("I126").Value = 5218.3335
$ws.Range("K126").Value = 15655.0005
$ws.Range("M126").Value = -10715.0005
$ws.Range("H128").Value = 427119.75
$ws.Range("I128").Value = 427119.75
$ws.Range("K128").Value = 1281359.25
$ws.Range("M128").Value = -1276379.25
$ws.Range("H135").Value = 994.4286
$ws.Range("I135").Value = 491.66666
$ws.Range("K135").Value = 4424.99994
$ws.Range("M135").Value = -1889.99994
$ws.Range("H137").Value = 1307.25
$ws.Range("I137").Value = 1208.3572
$ws.Range("J137").Value = 1999.5
$ws.Range("K137").Value = 3625.0716
$ws.Range("L137").Value = 5998.5
$ws.Range("M137").Value = 1474.9284
$ws.Range("N137").Value = -16198.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H103").Value = 25000
$ws.Range("J103").Value = 25000
$ws.Range("L103").Value = 25000
$ws.Range("N103").Value = -27344
$ws.Range("H123").Value = 42481.816
$ws.Range("J123").Value = 42481.816
$ws.Range("L123").Value = 42481.816
$ws.Range("N123").Value = -47381.816
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
$ws.Range("H132").Value = 7399.353
$ws.Range("I132").Value = 6863.636
$ws.Range("J132").Value = 8381.5
$ws.Range("K132").Value = 20590.908
$ws.Range("L132").Value = 25144.5
$ws.Range("M132").Value = -18060.908
$ws.Range("N132").Value = -30204.5
$ws.Range("H134").Value = 53997.6
$ws.Range("J134").Value = 53997.6
$ws.Range("L134").Value = 161992.8
$ws.Range("N134").Value = -167062.8
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 44999.9
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 44999.9
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 44999.9
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -45451.9
$ws.Range("H42").Value = 16644.334
$ws.Range("I42").Value = 17482
$ws.Range("J42").Value = 14969
$ws.Range("K42").Value = 17482
$ws.Range("L42").Value = 14969
$ws.Range("M42").Value = -16919
$ws.Range("N42").Value = -16095
$ws.Range("H49").Value = 16644.334
$ws.Range("I49").Value = 17482
$ws.Range("J49").Value = 14969
$ws.Range("K49").Value = 17482
$ws.Range("L49").Value = 14969
$ws.Range("M49").Value = -17335
$ws.Range("N49").Value = -15263
$ws.Range("H122").Value = 3479.8147
$ws.Range("J122").Value = 4600
$ws.Range("L122").Value = 13800
$ws.Range("N122").Value = -18700
$ws.Range("H132").Value = 10073
$ws.Range("I132").Value = 10175.4
$ws.Range("K132").Value = 30526.2
$ws.Range("M132").Value = -27996.2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 15751
$ws.Range("I81").Value = 11500
$ws.Range("K81").Value = 23000
$ws.Range("M81").Value = -21939
$ws.Range("H84").Value = 15751
$ws.Range("I84").Value = 11500
$ws.Range("K84").Value = 115000
$ws.Range("M84").Value = -109696
$ws.Range("H96").Value = 1462.8334
$ws.Range("J96").Value = 1375.6666
$ws.Range("L96").Value = 1375.6666
$ws.Range("N96").Value = -4121.6666
$ws.Range("H132").Value = 3796.7097
$ws.Range("I132").Value = 3627.0386
$ws.Range("K132").Value = 10881.1158
$ws.Range("M132").Value = -8351.1158
